$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HOURS")

# Fill in the IMPLEMENTATION row (row 8) hours that were forgotten on the 25th
$ws.Range("C8").Value = 165
$ws.Range("D8").Value = 205
$ws.Range("E8").Value = 70

# Move the active selection to E8, matching where editing left off
$ws.Range("E8").Select()
